# Append 13 more performance-result rows (V0 / unknown view) to the
# "Performance Results" sheet, right after the existing last row (507).
# Columns: A=(blank), B=IModel, C=View, D=Flags,
#          E=TileLoadingTime, F=Scene, G=GarbageExecute, H=InitCommands,
#          I=BackgroundDraw, J=SetClips, K=OpaqueDraw, L=TranslucentDraw,
#          M=HiliteDraw, N=CompositeDraw, O=OverlayDraw, P=RenderFrameTime,
#          Q=glFinish, R=TotalTime

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Performance Results")

$imodelPath = "d:\js\s\imodeljs-core/test-apps/testbed/frontend/performance/imodels/Wraith_MultiMulti.ibim"
$viewName = "V0"
$flags = "unknown"

$startRow = 508

$data = @(
    @(6144, 4, 0, 2, 0, 0, 0, 0, 13, 0, 0, 19, 0, 19),
    @(6144, 6, 0, 2, 1, 0, 0, 0, 15, 0, 0, 24, 0, 24),
    @(6144, 4, 0, 2, 0, 0, 0, 0, 14, 0, 0, 20, 0, 20),
    @(6144, 3, 0, 1, 0, 0, 0, 0, 12, 0, 0, 16, 0, 16),
    @(6144, 4, 0, 3, 0, 0, 0, 0, 104, 162, 0, 273, 184, 457),
    @(6144, 4, 0, 1, 1, 0, 0, 0, 22, 2, 0, 30, 0, 30),
    @(6144, 3, 0, 2, 0, 0, 0, 0, 13, 2, 0, 20, 0, 20),
    @(6144, 4, 0, 4, 1, 0, 0, 0, 17, 6, 0, 32, 0, 32),
    @(6144, 5, 0, 1, 0, 0, 0, 0, 18, 3, 0, 27, 0, 27),
    @(6144, 3, 0, 2, 0, 0, 0, 0, 21, 2, 0, 28, 0, 28),
    @(6144, 5, 0, 2, 0, 0, 0, 0, 14, 4, 0, 25, 0, 25),
    @(6144, 4, 0, 1, 0, 0, 0, 0, 12, 3, 0, 20, 0, 20),
    @(6144, 5, 0, 1, 0, 0, 0, 0, 14, 3, 0, 23, 0, 23)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Column A is blank text (matches the existing rows, which all store an
    # empty shared string rather than leaving the cell truly empty). A plain
    # `.Value = ""` assignment removes the cell instead of writing "", so use
    # the text-prefix apostrophe to force an empty text value, then reset the
    # style so it doesn't pick up the quote-prefix formatting flag.
    $ws.Cells.Item($r, 1).Value = "'"
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $imodelPath
    $ws.Cells.Item($r, 3).Value = $viewName
    $ws.Cells.Item($r, 4).Value = $flags

    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, 5 + $c).Value = $row[$c]
    }
}
